# "ext field driven flow"
# The sheet previously held one row of sample data (row 2) plus a stray
# value in A3. This clears that sample data out, leaving only the two
# cells that already carried a custom number format (E2 and H2) as empty,
# formatted placeholders, and extends the "H" (ext field) column with the
# same formatting down through row 16 so the sheet is ready to be driven
# by the external-field flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample row's values (A2:D2, F2:G2 incl. the G2 formula)
$ws.Range("A2:D2").ClearContents()
$ws.Range("F2:G2").ClearContents()

# E2 and H2 keep their existing (2-decimal) number format, just clear values
$ws.Range("E2").ClearContents()
$ws.Range("H2").ClearContents()

# Drop the leftover value in A3
$ws.Range("A3").ClearContents()

# Extend column H (same number format as H2) down through row 16
for ($r = 3; $r -le 16; $r++) {
    $ws.Range("H$r").NumberFormat = "0.00"
}

# Window/view bookkeeping to match the new layout
$win = $excel.ActiveWindow
$win.Width = 28800
$win.Height = 18000
$win.Top = 0
$win.Left = 0

# Move the active selection to C10
$ws.Range("C10").Select()
